$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (previously located right after
#    the "Số lượng hàng bảo quản" MERGEFIELD's fldChar end).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. In the "PHÒNG KỸ KTBQ" cell, drop the word "KỸ " so the text reads
#    "PHÒNG KTBQ", keeping the same bold/size/color formatting.
$rng = $d.Content.Duplicate
$rng.Find.Execute("PHÒNG KỸ KTBQ", $true, $false, $false, $false, $false, $true, 1, $false, "PHÒNG KTBQ", 2) | Out-Null

# 3. Re-insert a "_GoBack" bookmark between "PHÒNG " and "KTBQ", which
#    naturally splits the run into two runs around the bookmark.
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("PHÒNG ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPoint = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $insPoint) | Out-Null
